$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "Densidade"
$ws.Range("B2").Value = 1500
$ws.Range("C2").Value = 194.1122415647322
$ws.Range("D2").Value = 2000
$ws.Range("E2").Value = 1000
$ws.Range("G2").Value = 0.99
# H2 stays "Continuous"

# Delete rows 3 and 4 entirely
$ws.Rows("3:4").Delete()
